$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.999999946770628
$ws.Range("E2").Value = 0.999999946770628

$ws.Range("D3").Value = 0.9999999905803297
$ws.Range("E3").Value = 0.9999999905803297

$ws.Range("D4").Value = 0.0001370812719789731
$ws.Range("E4").Value = 0.0001370812719789731

$ws.Range("D5").Value = 0.000005119177639017285
$ws.Range("E5").Value = 0.000005119177639017285

$ws.Range("D6").Value = 0.001807892321715646
$ws.Range("E6").Value = 0.001807892321715646

$ws.Range("F7").Value = 5.871845245361328
